# Habitat sheet: add a new "user" column right after habitat_id (at B),
# pushing habitat_specific_type / habitat_main_type_id / approved one
# column to the right, and drop the old trailing "user" column by moving
# its data into the freshly inserted column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before column B. This shifts the existing
# B (habitat_specific_type), C (habitat_main_type_id) and D (approved)
# columns one place to the right (-> C, D, E) while keeping their custom
# widths, and also shifts the trailing "user" column from E to F.
$ws.Columns("B").Insert()

# Move the "user" column (now F, values + header) into the new column B.
$ws.Range("F1:F15").Cut()
$ws.Range("B1").Select()
$ws.Paste()

# Match the resulting selection: the whole of column B selected.
$ws.Columns("B").Select()
